# Edit: insert two new daily price rows (Damasco, variety "Dina") above the
# existing data block, pushing the previous rows 228-288 down to 230-290.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 228 (existing rows 228.. shift down to 230..)
$ws.Range("A228:A229").EntireRow.Insert()

# Common ("template") values shared by every data row in this sheet
$mercadoId = 6
$mercado   = "Mercado Mayorista Lo Valledor de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$tipo      = "Fruta"
$productoId = 100103
$producto   = "Frutos de hueso (carozo)"
$categoriaId = 100103003
$categoria   = "Damasco"

# New row 228
$r = 228
$ws.Cells.Item($r,1).Value  = $mercadoId
$ws.Cells.Item($r,2).Value  = $mercado
$ws.Cells.Item($r,3).Value  = $region
$ws.Cells.Item($r,4).Value  = 45258
$ws.Cells.Item($r,5).Value  = $codreg
$ws.Cells.Item($r,6).Value  = $tipo
$ws.Cells.Item($r,7).Value  = $productoId
$ws.Cells.Item($r,8).Value  = $producto
$ws.Cells.Item($r,9).Value  = $categoriaId
$ws.Cells.Item($r,10).Value = $categoria
$ws.Cells.Item($r,11).Value = "Dina"
$ws.Cells.Item($r,12).Value = "Especial"
$ws.Cells.Item($r,13).Value = 420
$ws.Cells.Item($r,14).Value = 25000
$ws.Cells.Item($r,15).Value = 25000
$ws.Cells.Item($r,16).Value = 25000
$ws.Cells.Item($r,17).Value = "$/caja 12 kilos"
$ws.Cells.Item($r,18).Value = "Región de O'Higgins"
$ws.Cells.Item($r,19).Value = 2083
$ws.Cells.Item($r,20).Value = 12

# New row 229
$r = 229
$ws.Cells.Item($r,1).Value  = $mercadoId
$ws.Cells.Item($r,2).Value  = $mercado
$ws.Cells.Item($r,3).Value  = $region
$ws.Cells.Item($r,4).Value  = 45258
$ws.Cells.Item($r,5).Value  = $codreg
$ws.Cells.Item($r,6).Value  = $tipo
$ws.Cells.Item($r,7).Value  = $productoId
$ws.Cells.Item($r,8).Value  = $producto
$ws.Cells.Item($r,9).Value  = $categoriaId
$ws.Cells.Item($r,10).Value = $categoria
$ws.Cells.Item($r,11).Value = "Dina"
$ws.Cells.Item($r,12).Value = "Primera"
$ws.Cells.Item($r,13).Value = 620
$ws.Cells.Item($r,14).Value = 23000
$ws.Cells.Item($r,15).Value = 23000
$ws.Cells.Item($r,16).Value = 23000
$ws.Cells.Item($r,17).Value = "$/caja 12 kilos"
$ws.Cells.Item($r,18).Value = "Región de O'Higgins"
$ws.Cells.Item($r,19).Value = 1917
$ws.Cells.Item($r,20).Value = 12

# Make sure the date cells keep the same date/time number format as the rest
# of column D.
$ws.Range("D228:D229").NumberFormat = "YYYY-MM-DD HH:MM:SS"
